# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
# Reordena la columna "Periodo Mora" (E16:E27) de orden descendente
# (2001,1912,...,1902) a orden ascendente (1902,1903,...,1912,2001).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periodos = @("1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912","2001")

$startRow = 16
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("E$row").Value = $periodos[$i]
}
